$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wellington's row (row 9) is missing his "Nota2" grade (column E).
# Fill it in with "5.5" (same text value/format used by the other students'
# grade cells, e.g. E4) to match the formatting of the other grade cells.
$ws.Range("E4").Copy($ws.Range("E9"))

# Update the active selection to reflect where the edit was made.
$ws.Range("E9").Select()
